$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: simple text fix inside "Requisito Funcional 6"
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "O jogo deve permitir ao jogador  ver os status e itens.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "O jogo deve permitir ao jogador  ver os status do seu personagem.",
    2)

# ---------------------------------------------------------------------------
# Change 2: flesh out "Requisito Funcional 8 / 9 / 10" with their body text.
#
# Each of these paragraphs currently looks like:
#   [bold] "Requisito Funcional N:  "   (two trailing spaces)
#   [normal] " "                        (single orphan space)
#
# For requisitos 8 and 9 we simply:
#   1. trim the bold label back to a single trailing space, and
#   2. overwrite the orphan trailing space with the real body text.
# ---------------------------------------------------------------------------

function Set-RequisitoBody($paraIndex, $label, $bodyText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $null = $r.Find.Execute(
        "$label`:  ", $true, $false, $false, $false, $false, $true, 1, $false,
        "$label`: ", 2)

    $p2 = $d.Paragraphs($paraIndex)
    $r2 = $p2.Range
    $tail = $d.Range($r2.End - 2, $r2.End - 1)
    $tail.Text = $bodyText
}

Set-RequisitoBody 133 "Requisito Funcional 8" "O jogo deve permitir ao jogador visualizar seus itens.  "
Set-RequisitoBody 135 "Requisito Funcional 9" "O jogo deve permitir ao jogador visualizar seus equipamentos.  "

# ---------------------------------------------------------------------------
# Requisito Funcional 10 ends up with an extra bold " " run between the new
# body text and the original trailing space run, so it needs one more step.
# ---------------------------------------------------------------------------

$p = $d.Paragraphs(137)
$r = $p.Range
$null = $r.Find.Execute(
    "Requisito Funcional 10:  ", $true, $false, $false, $false, $false, $true, 1, $false,
    "Requisito Funcional 10: ", 2)

$p = $d.Paragraphs(137)
$r = $p.Range
$boundary = $r.Start + "Requisito Funcional 10: ".Length

$newChunk = "O jogo deve permitir ao jogador visualizar suas habilidades.  "
$insertPoint = $d.Range($boundary, $boundary)
$insertPoint.Text = $newChunk

$bodyLen = $newChunk.Length - 1
$bodyRange = $d.Range($boundary, $boundary + $bodyLen)
$bodyRange.Bold = 0
